$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new product rows at 144:145 (list is sorted alphabetically, these two
#     new items sort in right before the existing "كالونا" row) ---
$ws.Rows("144:145").Insert()

# Copy cell formatting (styles/borders/number formats) from the row that used to be 144
# (now shifted to 146) into the two freshly inserted blank rows so they look identical
# to every other data row in the table.
$ws.Range("A146:Q146").Copy()
$ws.Range("A144").PasteSpecial(-4122)
$ws.Range("A145").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the row heights used by the new rows.
$ws.Rows("144").RowHeight = 24.75
$ws.Rows("145").RowHeight = 25.5

# Recreate the merged-cell layout (A:B, C:G, H:K, L:M, N:O) used by every data row.
foreach ($r in 144,145) {
  $ws.Range("A" + $r + ":B" + $r).Merge()
  $ws.Range("C" + $r + ":G" + $r).Merge()
  $ws.Range("H" + $r + ":K" + $r).Merge()
  $ws.Range("L" + $r + ":M" + $r).Merge()
  $ws.Range("N" + $r + ":O" + $r).Merge()
}

# --- Fill in the values for the two new rows ---
# Row 144: "صوفى طويل جدا جدا"
$ws.Range("A144").Value = 138
$ws.Range("C144").Value = "صوفى طويل جدا جدا"
$ws.Range("H144").Value = "11:0"
$ws.Range("L144").NumberFormat = "@"
$ws.Range("L144").Value = "0"
$ws.Range("L144").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N144").Value = "55.00"
$ws.Range("P144").NumberFormat = "@"
$ws.Range("P144").Value = "55.0000"
$ws.Range("P144").NumberFormat = "0.00"
$ws.Range("Q144").Value = "1:0"

# Row 145: "فلامنجو شفرات للنساء"
$ws.Range("A145").Value = 139
$ws.Range("C145").Value = "فلامنجو شفرات للنساء"
$ws.Range("H145").Value = "5:0"
$ws.Range("L145").NumberFormat = "@"
$ws.Range("L145").Value = "0"
$ws.Range("L145").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"
$ws.Range("N145").Value = "40.00"
$ws.Range("P145").NumberFormat = "@"
$ws.Range("P145").Value = "40.0000"
$ws.Range("P145").NumberFormat = "0.00"
$ws.Range("Q145").Value = "1:0"

# --- Update the grand-total cell (shifted from row 150 to row 152) ---
$ws.Range("P152").Value = 9950.83

# --- Update the printed timestamp in the footer (shifted from row 151 to row 153) ---
$ws.Range("A153").Value = "Monday, 11 August, 2025 10:25 PM"
